$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 & 11 (Status column F): "in Unity but incomplete" -> "implemented in Unity"
$ws.Range("F10").Value = "implemented in Unity"
$ws.Range("F11").Value = "implemented in Unity"

# Row 15 (Description column C): "small bump noise" -> "small bump noise with rock"
$ws.Range("C15").Value = "small bump noise with rock"

# Update the active selection to C12
$ws.Range("C12").Select()
